$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (number wanting to go) in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1230
$ws1.Range("F17").Value = 1072
$ws1.Range("F18").Value = 4570
$ws1.Range("F23").Value = 436
$ws1.Range("F24").Value = 1202
$ws1.Range("F26").Value = 2575
$ws1.Range("F29").Value = 104
$ws1.Range("F37").Value = 2496
$ws1.Range("F38").Value = 2247

# Sheet "全部类型" (All types) - same rows (same events, shifted by one row)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1230
$ws4.Range("F18").Value = 1072
$ws4.Range("F19").Value = 4570
$ws4.Range("F26").Value = 436
$ws4.Range("F27").Value = 1202
$ws4.Range("F29").Value = 2575
$ws4.Range("F32").Value = 104
$ws4.Range("F42").Value = 2496
$ws4.Range("F44").Value = 2247
